$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add I1 = "I0" and J1 = "IF" ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold, centered, bordered style) from the existing
# header cell H1 onto the two new header cells so they match the rest
# of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data rows (rows 2-31): add columns I (I0) and J (IF) ---
$iValues = @(6,7,4,4,4,7,4,5,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,4,1)
$jValues = @(7,8,5,7,6,7,6,6,3,3,5,4,4,5,2,5,6,6,3,5,7,7,7,7,5,3,5,5,6,2)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
